# Scheduled runner update: refresh profit calculations across several
# sheets (market-board price pulls changed, so downstream profit math
# shifts) and backfill newly-priced leves that previously had blank
# price/profit columns.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        $ws,
        [int]$row,
        [hashtable]$values
    )
    foreach ($col in $values.Keys) {
        $ws.Cells.Item($row, $col).Value = $values[$col]
    }
}

# Column indices: H=8 I=9 J=10 K=11 L=12 M=13 N=14

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")

Set-Row $ws 74 @{ 8=3563.913; 9=2993.889; 10=5616; 11=2993.889; 12=5616; 13=-2057.889; 14=-7488 }
Set-Row $ws 77 @{ 8=3563.913; 9=2993.889; 10=5616; 11=14969.445; 12=28080; 13=-10289.445; 14=-37440 }
Set-Row $ws 86 @{ 8=9060.691999999999; 9=1383.3334; 10=15641.286; 11=1383.3334; 12=15641.286; 13=-260.3334; 14=-17887.286 }
Set-Row $ws 89 @{ 8=9060.691999999999; 9=1383.3334; 10=15641.286; 11=6916.666999999999; 12=78206.42999999999; 13=-1300.666999999999; 14=-89438.42999999999 }
Set-Row $ws 137 @{ 8=25327848; 9=8334221; 10=46569884; 11=25002663; 12=139709652; 13=-25000113; 14=-139714752 }

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")

Set-Row $ws 32 @{ 8=1059.6; 9=1059.6; 11=1059.6; 13=-772.5999999999999 }
Set-Row $ws 88 @{ 8=3835.077; 9=2000; 10=4385.6; 11=2000; 12=4385.6; 13=-1594; 14=-5197.6 }
Set-Row $ws 91 @{ 8=3835.077; 9=2000; 10=4385.6; 11=2000; 12=4385.6; 13=-596; 14=-7193.6 }
Set-Row $ws 102 @{ 8=2818.8823; 9=2574.1667; 10=3406.2; 11=2574.1667; 12=3406.2; 13=-952.1667000000002; 14=-6650.2 }

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")

Set-Row $ws 86 @{ 8=1955.05; 9=1961.2858; 10=1649.5; 11=1961.2858; 12=1649.5; 13=-838.2858000000001; 14=-3895.5 }
Set-Row $ws 89 @{ 8=1955.05; 9=1961.2858; 10=1649.5; 11=9806.429; 12=8247.5; 13=-4190.429; 14=-19479.5 }
Set-Row $ws 134 @{ 8=11804970; 9=12656163; 10=172000; 11=37968489; 12=516000; 13=-37965954; 14=-521070 }

# ---------------- CUL ----------------
$ws = $wb.Worksheets.Item("CUL")

Set-Row $ws 131 @{ 8=5113251.5; 9=62625350; 10=1064.6; 11=187876050; 12=3193.8; 13=-187871010; 14=-13273.8 }
Set-Row $ws 132 @{ 8=1728.2858; 9=980; 10=2027.6; 11=8820; 12=18248.4; 13=-6290; 14=-23308.4 }

# ---------------- GSM (newly backfilled leves, rows 125-141) ----------------
$ws = $wb.Worksheets.Item("GSM")

Set-Row $ws 125 @{ 8=35686.8; 9=0; 10=35686.8; 11=0; 12=35686.8; 14=-40606.8 }
Set-Row $ws 126 @{ 8=6666.375; 9=11444.9; 10=3253.1428; 11=34334.7; 12=9759.4284; 13=-31864.7; 14=-14699.4284 }
Set-Row $ws 127 @{ 8=54980; 9=0; 10=54980; 11=0; 12=54980; 14=-64900 }
Set-Row $ws 128 @{ 8=41280; 9=0; 10=41280; 11=0; 12=41280; 14=-51240 }
Set-Row $ws 129 @{ 8=34254; 9=0; 10=34254; 11=0; 12=34254; 14=-44254 }
Set-Row $ws 130 @{ 8=49980; 9=0; 10=49980; 11=0; 12=49980; 14=-60020 }
Set-Row $ws 131 @{ 8=0; 9=0; 10=0; 11=0; 12=0 }
Set-Row $ws 132 @{ 8=5884391.5; 9=9525537; 10=2541.8462; 11=28576611; 12=7625.5386; 13=-28574081; 14=-12685.5386 }
Set-Row $ws 133 @{ 8=51167.8; 9=0; 10=51167.8; 11=0; 12=51167.8; 14=-61287.8 }
Set-Row $ws 134 @{ 8=24719.9; 9=0; 10=24719.9; 11=0; 12=74159.70000000001; 14=-79229.70000000001 }
Set-Row $ws 135 @{ 8=71000; 9=0; 10=71000; 11=0; 12=71000; 14=-81140 }
Set-Row $ws 136 @{ 8=22946.525; 9=0; 10=22946.525; 11=0; 12=68839.57500000001; 14=-73939.57500000001 }
Set-Row $ws 137 @{ 8=0; 9=0; 10=0; 11=0; 12=0 }
Set-Row $ws 138 @{ 8=39345.105; 9=0; 10=39345.105; 11=0; 12=39345.105; 14=-49625.105 }
Set-Row $ws 139 @{ 8=67021.75; 9=300; 10=76553.42999999999; 11=300; 12=76553.42999999999; 13=4840; 14=-86833.42999999999 }
Set-Row $ws 140 @{ 8=67100; 9=0; 10=67100; 11=0; 12=67100; 14=-77460 }
Set-Row $ws 141 @{ 8=51300; 9=0; 10=51300; 11=0; 12=51300; 14=-61660 }

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")

Set-Row $ws 16 @{ 8=35286.723; 9=50572.6; 11=50572.6; 13=-50402.6 }

Write-Host "Profit sheets refreshed."
